$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.357.52'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.695.39'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.06'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5421'
$ws.Range('E6').Value = '  +2.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.009'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2735'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06461'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.67'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07667'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.723.80'
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.538'
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5818'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008418'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.10'
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.402.99'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.923'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.009'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.90'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.04'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.285'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.009'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.04'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1288'
$ws.Range('E25').Value = '  +3.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.856'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.90'
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06354'
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.387'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.327'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.616'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.599'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.691'
$ws.Range('E33').Value = '  +1.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.034'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6204'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.417'
$ws.Range('E36').Value = '  +0.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.755'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01653'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.116.21'
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.103'
$ws.Range('E40').Value = '  -5.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8875'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.014'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.17'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.847.90'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.79'
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.200'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05287'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.103'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4302'
$ws.Range('E51').Value = '  +0.22%  '
